# 2021FSAdates.xlsx - "Add files via upload" commit
#
# The underlying Table1 (xml-mapped) data was edited on the "Sheet1"
# worksheet:
#   - Row 93  (20211107 U13 E): SplitGender (E93) changed 1 -> 0
#     (so the event is no longer split by gender; K/L helper formulas
#     recompute automatically).
#   - Rows 103, 104, 105, 107 (20211121 OMS/VMF/OMF/VMS): Cancelled
#     (D<row>) changed 0 -> 1 (these events were cancelled; K/L helper
#     formulas recompute automatically to "N/A"/"").
#   - The sheet's view scrolled down a bit and the active selection moved
#     from K101 to E106.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Data edits -----------------------------------------------------

# Row 93: turn off SplitGender for the 20211107 U13 Epee event.
$ws.Range("E93").Value = 0

# Rows 103, 104, 105, 107: mark these 20211121 events as Cancelled.
$ws.Range("D103").Value = 1
$ws.Range("D104").Value = 1
$ws.Range("D105").Value = 1
$ws.Range("D107").Value = 1

# --- View / selection -------------------------------------------------

# Scroll the window so row 82 is at the top and select E106 (matches the
# author's on-screen state when the file was saved).
try {
    $excel.ActiveWindow.ScrollRow = 82
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

[void]$ws.Range("E106").Select()
